$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.011.19'
$ws.Range('E2').Value = '  +1.41%  '
$ws.Range('D3').Value = '1.888.74'
$ws.Range('E3').Value = '  +1.11%  '
$ws.Range('D4').Formula = '="0.9995"'
$ws.Range('D4').Copy()
$ws.Range('D4').PasteSpecial(-4163)
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Formula = '="331.37"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  -1.98%  '
$ws.Range('D6').Formula = '="0.9993"'
$ws.Range('D6').Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').Formula = '="0.4604"'
$ws.Range('D7').Copy()
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('E7').Value = '  -1.72%  '
$ws.Range('D8').Formula = '="0.4099"'
$ws.Range('D8').Copy()
$ws.Range('D8').PasteSpecial(-4163)
$ws.Range('E8').Value = '  +2.79%  '
$ws.Range('D9').Formula = '="47.35"'
$ws.Range('D9').Copy()
$ws.Range('D9').PasteSpecial(-4163)
$ws.Range('E9').Value = '  -0.85%  '
$ws.Range('D10').Formula = '="0.07988"'
$ws.Range('D10').Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('E10').Value = '  -0.56%  '
$ws.Range('D11').Formula = '="0.9908"'
$ws.Range('D11').Copy()
$ws.Range('D11').PasteSpecial(-4163)
$ws.Range('E11').Value = '  -0.96%  '
$ws.Range('D12').Formula = '="21.72"'
$ws.Range('D12').Copy()
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('E12').Value = '  -1.70%  '
$ws.Range('D13').Value = '1.883.67'
$ws.Range('E13').Value = '  +1.32%  '
$ws.Range('D14').Formula = '="5.907"'
$ws.Range('D14').Copy()
$ws.Range('D14').PasteSpecial(-4163)
$ws.Range('E14').Value = '  -2.31%  '
$ws.Range('D15').Formula = '="7.077"'
$ws.Range('D15').Copy()
$ws.Range('D15').PasteSpecial(-4163)
$ws.Range('E15').Value = '  -2.90%  '
$ws.Range('E16').Value = '  -1.38%  '
$ws.Range('D17').Formula = '="1.000"'
$ws.Range('D17').Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('D18').Formula = '="0.00001029"'
$ws.Range('D18').Copy()
$ws.Range('D18').PasteSpecial(-4163)
$ws.Range('E18').Value = '  -1.25%  '
$ws.Range('D19').Formula = '="0.06557"'
$ws.Range('D19').Copy()
$ws.Range('D19').PasteSpecial(-4163)
$ws.Range('E19').Value = '  -0.82%  '
$ws.Range('E20').Value = '  -0.75%  '
$ws.Range('D21').Formula = '="1.000"'
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('D22').Value = '29.040.35'
$ws.Range('E22').Value = '  +1.43%  '
$ws.Range('D23').Formula = '="5.401"'
$ws.Range('D23').Copy()
$ws.Range('D23').PasteSpecial(-4163)
$ws.Range('E23').Value = '  -1.59%  '
$ws.Range('D24').Formula = '="11.23"'
$ws.Range('D24').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  +1.82%  '
$ws.Range('D25').Formula = '="2.212"'
$ws.Range('D25').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  -1.82%  '
$ws.Range('D26').Value = '2.121.11'
$ws.Range('E26').Value = '  +1.85%  '
$ws.Range('D27').Formula = '="157.15"'
$ws.Range('D27').Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E27').Value = '  -2.29%  '
$ws.Range('E28').Value = '  -0.63%  '
$ws.Range('D29').Formula = '="2.106"'
$ws.Range('D29').Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  -0.54%  '
$ws.Range('E30').Value = '  -1.52%  '
$ws.Range('D31').Formula = '="117.96"'
$ws.Range('D31').Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('E31').Value = '  -1.44%  '
$ws.Range('D32').Formula = '="0.9760"'
$ws.Range('D32').Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('E32').Value = '  +0.33%  '
$ws.Range('D33').Formula = '="0.09342"'
$ws.Range('D33').Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('E33').Value = '  -1.91%  '
$ws.Range('E34').Value = '  -0.27%  '
$ws.Range('D35').Formula = '="1.410"'
$ws.Range('D35').Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('E35').Value = '  +2.04%  '
$ws.Range('D36').Formula = '="5.277"'
$ws.Range('D36').Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('E36').Value = '  -1.61%  '
$ws.Range('D37').Formula = '="0.06047"'
$ws.Range('D37').Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Value = '  -2.51%  '
$ws.Range('D38').Formula = '="0.02226"'
$ws.Range('D38').Copy()
$ws.Range('D38').PasteSpecial(-4163)
$ws.Range('E38').Value = '  -1.24%  '
$ws.Range('D39').Formula = '="8.244"'
$ws.Range('D39').Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Value = '  -1.88%  '
$ws.Range('D40').Formula = '="1.181"'
$ws.Range('D40').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('D41').Formula = '="0.9984"'
$ws.Range('D41').Copy()
$ws.Range('D41').PasteSpecial(-4163)
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('D42').Formula = '="0.5773"'
$ws.Range('D42').Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('E42').Value = '  -2.73%  '
$ws.Range('D43').Formula = '="10.14"'
$ws.Range('D43').Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('E43').Value = '  -1.75%  '
$ws.Range('D44').Formula = '="0.1823"'
$ws.Range('D44').Copy()
$ws.Range('D44').PasteSpecial(-4163)
$ws.Range('E44').Value = '  -3.13%  '
$ws.Range('D45').Formula = '="1.258"'
$ws.Range('D45').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('D46').Formula = '="2.274"'
$ws.Range('D46').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  +9.03%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Formula = '="0.5465"'
$ws.Range('D47').Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('E47').Value = '  -1.66%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Formula = '="11.98"'
$ws.Range('D48').Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Value = '  -1.73%  '
$ws.Range('D49').Formula = '="1.896"'
$ws.Range('D49').Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('D50').Formula = '="0.07029"'
$ws.Range('D50').Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('E50').Value = '  -5.56%  '
$ws.Range('D51').Formula = '="45.35"'
$ws.Range('D51').Copy()
$ws.Range('D51').PasteSpecial(-4163)
$ws.Range('E51').Value = '  +13.93%  '

$excel.CutCopyMode = 0
